$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column widths (col A and col B).
# The host rounds ColumnWidth to the nearest pixel, so the inputs below are
# chosen to land as close as possible to the target widths
# (15.42578125 and 14.7109375 "characters" respectively).
$ws.Columns(1).ColumnWidth = 14.666666666666668
$ws.Columns(2).ColumnWidth = 13.833333333333332

# Update the computed values in columns A and B
$ws.Range("A1").Value2 = -0.094724588703769541
$ws.Range("B1").Value2 = 0.094689250771075706
$ws.Range("A2").Value2 = -0.072584806478448449
$ws.Range("B2").Value2 = 0.072497166465987917
$ws.Range("A3").Value2 = -0.022795891204763663
$ws.Range("B3").Value2 = 0.022708145759018805
$ws.Range("A4").Value2 = -0.014708145806235251
$ws.Range("B4").Value2 = 0.014222713428843292
$ws.Range("A5").Value2 = -0.011222713449978272
$ws.Range("B5").Value2 = 0.0095595318596117806
$ws.Range("A6").Value2 = 0.0042136224721893711
$ws.Range("B6").Value2 = -0.0044441287699168441
$ws.Range("A7").Value2 = 0.014444128709740767
$ws.Range("B7").Value2 = -0.014494364527986914
$ws.Range("A8").Value2 = 0.024494364468697682
$ws.Range("B8").Value2 = -0.024584211120325605
$ws.Range("A9").Value2 = 0.026584211106190025
$ws.Range("B9").Value2 = -0.026660168386907657
$ws.Range("A10").Value2 = -0.027394431567280719
$ws.Range("B10").Value2 = 0.027379426963408449
$ws.Range("A11").Value2 = -0.024379426980860153
$ws.Range("B11").Value2 = 0.024353951516760119
$ws.Range("A12").Value2 = -0.020853951537453064
$ws.Range("B12").Value2 = 0.020665557789395272
$ws.Range("A13").Value2 = -0.017165557813531684
$ws.Range("B13").Value2 = 0.017079780331878069
$ws.Range("A14").Value2 = -0.0090797803811168976
$ws.Range("B14").Value2 = 0.0090517884760918932
$ws.Range("A15").Value2 = -0.0080517884883324342
$ws.Range("B15").Value2 = 0.0080337141027522208
$ws.Range("A16").Value2 = -0.0060337141209161338
$ws.Range("B16").Value2 = 0.0060034120637779687
$ws.Range("A17").Value2 = -0.0040034120827785458
$ws.Range("B17").Value2 = 0.003999999970089263
$ws.Range("A18").Value2 = -0.016104444324842149
$ws.Range("B18").Value2 = 0.016091306863600607
$ws.Range("A19").Value2 = -0.012091306885170017
$ws.Range("B19").Value2 = 0.012016646345022419
$ws.Range("A20").Value2 = -0.0080166463684534506
$ws.Range("B20").Value2 = 0.0080056593283686084
$ws.Range("A21").Value2 = -0.0040056593520851891
$ws.Range("B21").Value2 = 0.0039999999761084482
$ws.Range("A22").Value2 = -0.045701275286120335
$ws.Range("B22").Value2 = 0.045491205525898692
$ws.Range("A23").Value2 = -0.040491205557957599
$ws.Range("B23").Value2 = 0.040097460105165794
$ws.Range("A24").Value2 = -0.020097460221417229
$ws.Range("B24").Value2 = 0.019999999882048591
$ws.Range("A25").Value2 = -0.011773154349036474
$ws.Range("B25").Value2 = 0.011695125553067953
$ws.Range("A26").Value2 = -0.0091951255731608228
$ws.Range("B26").Value2 = 0.0090969649923913209
$ws.Range("A27").Value2 = -0.0065969650127577495
$ws.Range("B27").Value2 = 0.0060297960836739328
$ws.Range("A28").Value2 = -0.0040297961023636475
$ws.Range("B28").Value2 = 0.0036567513113707761
$ws.Range("A29").Value2 = -0.056675798260207522
$ws.Range("B29").Value2 = 0.056544884846355892
$ws.Range("A30").Value2 = -0.021165912848355006
$ws.Range("B30").Value2 = 0.021022693030058992
$ws.Range("A31").Value2 = -0.014022693081180648
$ws.Range("B31").Value2 = 0.014001239496391094
$ws.Range("A32").Value2 = -0.004001239563999448
$ws.Range("B32").Value2 = 0.0039999999644955153
